$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: only the scrape timestamp changed ---
$ws.Range('O2').Value = '2022-08-20 20:58:05'

# --- Row 3: product data refreshed by the crawler ---
$ws.Range('A3').Value = '''6031467017'
$ws.Range('B3').Value = 'Naturaline Herren T-Shirt Kurzarm schwarz M'
$ws.Range('C3').Value = '/de/haushalt-tier/bekleidung/shirts-pullover/herren-shirt/naturaline-herren-t-shirt-kurzarm-schwarz-m/p/6031467017'
$ws.Range('D3').ClearContents() | Out-Null
$ws.Range('E3').ClearContents() | Out-Null
$ws.Range('F3').Value = 0
$ws.Range('G3').Value = 'Coop'
$ws.Range('H3').Value = '''24.95'
$ws.Range('I3').ClearContents() | Out-Null
$ws.Range('J3').ClearContents() | Out-Null
$ws.Range('K3').ClearContents() | Out-Null
$ws.Range('L3').ClearContents() | Out-Null
$ws.Range('M3').Value = '[''haushalt-tier'', ''bekleidung'', ''shirts-pullover'', ''herren-shirt'']'
$ws.Range('N3').Value = 'Naturaline Herren T-Shirt Kurzarm schwarz M - Online kein Bestand 24.95 Schweizer Franken'
$ws.Range('O3').Value = '2022-08-20 20:58:05'

# --- Row 4: product data refreshed by the crawler ---
$ws.Range('A4').Value = '''6077154005'
$ws.Range('B4').Value = 'Avela Söckchen Pure Noir One Size'
$ws.Range('C4').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/socken/avela-soeckchen-pure-noir-one-size/p/6077154005'
$ws.Range('D4').Value = '2ST'
$ws.Range('E4').ClearContents() | Out-Null
$ws.Range('F4').Value = 0
$ws.Range('G4').Value = 'Avela'
$ws.Range('H4').Value = '''2.95'
$ws.Range('I4').Value = '1.48/1ST'
$ws.Range('J4').Value = 'Preis pro 1 Stück'
$ws.Range('K4').Value = '''1.48'
$ws.Range('L4').Value = '1ST'
$ws.Range('M4').Value = '[''haushalt-tier'', ''bekleidung'', ''socken-unterwaesche'', ''socken'']'
$ws.Range('N4').Value = 'Avela Söckchen Pure Noir One Size 2.95 Schweizer Franken'
$ws.Range('O4').Value = '2022-08-20 20:58:05'

# --- Row 5: product data refreshed by the crawler ---
$ws.Range('A5').Value = '''6075745012'
$ws.Range('B5').Value = 'Avela Strumpfhose Madame Natural  11 - 12'
$ws.Range('C5').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-madame-natural-11-12/p/6075745012'
$ws.Range('D5').ClearContents() | Out-Null
$ws.Range('E5').ClearContents() | Out-Null
$ws.Range('F5').Value = 0
$ws.Range('G5').Value = 'Avela'
$ws.Range('H5').Value = '''5.95'
$ws.Range('I5').ClearContents() | Out-Null
$ws.Range('J5').ClearContents() | Out-Null
$ws.Range('K5').ClearContents() | Out-Null
$ws.Range('L5').ClearContents() | Out-Null
$ws.Range('M5').Value = '[''haushalt-tier'', ''bekleidung'', ''socken-unterwaesche'', ''struempfe'']'
$ws.Range('N5').Value = 'Avela Strumpfhose Madame Natural  11 - 12 5.95 Schweizer Franken'
$ws.Range('O5').Value = '2022-08-20 20:58:05'

# --- Row 6: product data refreshed by the crawler ---
$ws.Range('A6').Value = '''5799901001'
$ws.Range('B6').Value = 'Magic Matic navy ecorepel'
$ws.Range('C6').Value = '/de/haushalt-tier/bekleidung/taschen-accessoires/schirme/magic-matic-navy-ecorepel/p/5799901001'
$ws.Range('D6').ClearContents() | Out-Null
$ws.Range('E6').Value = 2
$ws.Range('F6').Value = 5
$ws.Range('G6').Value = 'Coop'
$ws.Range('H6').Value = '''19.95'
$ws.Range('I6').ClearContents() | Out-Null
$ws.Range('J6').ClearContents() | Out-Null
$ws.Range('K6').ClearContents() | Out-Null
$ws.Range('L6').ClearContents() | Out-Null
$ws.Range('M6').Value = '[''haushalt-tier'', ''bekleidung'', ''taschen-accessoires'', ''schirme'']'
$ws.Range('N6').Value = 'Magic Matic navy ecorepel 19.95 Schweizer Franken'
$ws.Range('O6').Value = '2022-08-20 20:58:05'

# --- Rows 7-73: only the scrape timestamp column changed ---
for ($r = 7; $r -le 73; $r++) {
    $ws.Cells.Item($r, 15).Value = '2022-08-20 20:58:05'
}

Write-Output "edit applied"